$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "E-T1-G101",
    "E-T1-G102",
    "E-T2-G103",
    "E-T2-G104",
    "E-T3-G105",
    "E-T3-G106",
    "E-T4-G107",
    "E-T4-G108",
    "E-T5-G109",
    "E-T5-G110"
)

$row = 101
foreach ($val in $values) {
    $ws.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}
